$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X23").Value = 1040
$ws.Range("AB23").Value = 3603
$ws.Range("X24").Value = 1305
$ws.Range("AB24").Value = 4775
$ws.Range("X59").Value = 5146
$ws.Range("AB59").Value = 28223
$ws.Range("X60").Value = 5180
$ws.Range("AB60").Value = 28406
$ws.Range("X61").Value = 5201
$ws.Range("AB61").Value = 28569
$ws.Range("X62").Value = 5212
$ws.Range("AB62").Value = 28699
$ws.Range("X63").Value = 5212
$ws.Range("AB63").Value = 28762
$ws.Range("AB64").Value = 28797

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("X36").Value = 93
$ws.Range("AB36").Value = 436
$ws.Range("X37").Value = 105
$ws.Range("AB37").Value = 502
$ws.Range("X38").Value = 116
$ws.Range("AB38").Value = 565
$ws.Range("X39").Value = 135
$ws.Range("AB39").Value = 626
$ws.Range("X40").Value = 153
$ws.Range("AB40").Value = 687
$ws.Range("X41").Value = 164
$ws.Range("AB41").Value = 752
$ws.Range("X42").Value = 176
$ws.Range("AB42").Value = 804
$ws.Range("X43").Value = 192
$ws.Range("AB43").Value = 862
$ws.Range("X44").Value = 204
$ws.Range("AB44").Value = 922
$ws.Range("X45").Value = 216
$ws.Range("AB45").Value = 989
$ws.Range("X46").Value = 230
$ws.Range("AB46").Value = 1040
$ws.Range("X47").Value = 237
$ws.Range("AB47").Value = 1089
$ws.Range("X48").Value = 243
$ws.Range("AB48").Value = 1126
$ws.Range("X49").Value = 251
$ws.Range("AB49").Value = 1181
$ws.Range("X50").Value = 257
$ws.Range("AB50").Value = 1211
$ws.Range("X51").Value = 267
$ws.Range("AB51").Value = 1249
$ws.Range("X52").Value = 284
$ws.Range("AB52").Value = 1301
$ws.Range("X53").Value = 297
$ws.Range("AB53").Value = 1346
$ws.Range("X54").Value = 302
$ws.Range("AB54").Value = 1391
$ws.Range("X55").Value = 310
$ws.Range("AB55").Value = 1433
$ws.Range("X56").Value = 315
$ws.Range("AB56").Value = 1455
$ws.Range("X57").Value = 326
$ws.Range("AB57").Value = 1496
$ws.Range("X58").Value = 334
$ws.Range("AB58").Value = 1535
$ws.Range("X59").Value = 338
$ws.Range("AB59").Value = 1560
$ws.Range("X60").Value = 347
$ws.Range("AB60").Value = 1586
$ws.Range("X61").Value = 349
$ws.Range("AB61").Value = 1608
$ws.Range("X62").Value = 354
$ws.Range("AB62").Value = 1625
$ws.Range("X63").Value = 355
$ws.Range("AB63").Value = 1634
$ws.Range("AB64").Value = 1640

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X28").Value = 230
$ws.Range("AB28").Value = 1090
$ws.Range("X29").Value = 248
$ws.Range("AB29").Value = 1207
$ws.Range("X30").Value = 284
$ws.Range("AB30").Value = 1364
$ws.Range("X32").Value = 325
$ws.Range("AB32").Value = 1611
$ws.Range("X33").Value = 360
$ws.Range("AB33").Value = 1794
$ws.Range("X34").Value = 369
$ws.Range("AB34").Value = 1882
$ws.Range("X35").Value = 371
$ws.Range("AB35").Value = 1996
$ws.Range("X36").Value = 374
$ws.Range("AB36").Value = 2181
$ws.Range("X37").Value = 382
$ws.Range("AB37").Value = 2209
$ws.Range("X38").Value = 388
$ws.Range("AB38").Value = 2290
$ws.Range("X39").Value = 379
$ws.Range("AB39").Value = 2354
$ws.Range("X40").Value = 370
$ws.Range("AB40").Value = 2335
$ws.Range("X41").Value = 370
$ws.Range("AB41").Value = 2323
$ws.Range("X42").Value = 381
$ws.Range("AB42").Value = 2308
$ws.Range("X43").Value = 365
$ws.Range("AB43").Value = 2312
$ws.Range("X44").Value = 338
$ws.Range("AB44").Value = 2252
$ws.Range("X45").Value = 331
$ws.Range("AB45").Value = 2142
$ws.Range("X46").Value = 314
$ws.Range("AB46").Value = 2073
$ws.Range("X47").Value = 305
$ws.Range("AB47").Value = 2013
$ws.Range("X48").Value = 291
$ws.Range("AB48").Value = 1938
$ws.Range("X49").Value = 300
$ws.Range("AB49").Value = 1915
$ws.Range("X50").Value = 300
$ws.Range("AB50").Value = 1893
$ws.Range("X51").Value = 285
$ws.Range("AB51").Value = 1860
$ws.Range("X52").Value = 258
$ws.Range("AB52").Value = 1736
$ws.Range("X53").Value = 252
$ws.Range("AB53").Value = 1680
$ws.Range("X54").Value = 236
$ws.Range("AB54").Value = 1582
$ws.Range("X55").Value = 233
$ws.Range("AB55").Value = 1531
$ws.Range("X56").Value = 236
$ws.Range("AB56").Value = 1525
$ws.Range("X57").Value = 219
$ws.Range("AB57").Value = 1502
$ws.Range("X58").Value = 207
$ws.Range("AB58").Value = 1426
$ws.Range("X61").Value = 181
$ws.Range("AB61").Value = 1259
$ws.Range("X62").Value = 185
$ws.Range("AB62").Value = 1239
$ws.Range("X63").Value = 192
$ws.Range("AB63").Value = 1221
$ws.Range("AB64").Value = 1224

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X28").Value = 30
$ws.Range("AB28").Value = 146
$ws.Range("X29").Value = 39
$ws.Range("AB29").Value = 175
$ws.Range("X30").Value = 40
$ws.Range("AB30").Value = 178
$ws.Range("X31").Value = 51
$ws.Range("AB31").Value = 211
$ws.Range("X32").Value = 55
$ws.Range("AB32").Value = 233
$ws.Range("X33").Value = 59
$ws.Range("AB33").Value = 254
$ws.Range("X34").Value = 61
$ws.Range("AB34").Value = 281
$ws.Range("X35").Value = 62
$ws.Range("AB35").Value = 291
$ws.Range("X36").Value = 65
$ws.Range("AB36").Value = 338
$ws.Range("X37").Value = 68
$ws.Range("AB37").Value = 360
$ws.Range("X38").Value = 71
$ws.Range("AB38").Value = 387
$ws.Range("X39").Value = 72
$ws.Range("AB39").Value = 393
$ws.Range("X40").Value = 70
$ws.Range("AB40").Value = 399
$ws.Range("X41").Value = 70
$ws.Range("AB41").Value = 399
$ws.Range("X42").Value = 69
$ws.Range("AB42").Value = 389
$ws.Range("X43").Value = 67
$ws.Range("AB43").Value = 389
$ws.Range("X44").Value = 66
$ws.Range("AB44").Value = 381
$ws.Range("X45").Value = 65
$ws.Range("AB45").Value = 382
$ws.Range("X46").Value = 65
$ws.Range("AB46").Value = 375
$ws.Range("X47").Value = 67
$ws.Range("AB47").Value = 368
$ws.Range("X48").Value = 67
$ws.Range("AB48").Value = 368
$ws.Range("X49").Value = 66
$ws.Range("AB49").Value = 359
$ws.Range("X50").Value = 62
$ws.Range("AB50").Value = 354
$ws.Range("X59").Value = 46
$ws.Range("AB59").Value = 211
$ws.Range("X61").Value = 39
$ws.Range("AB61").Value = 198
$ws.Range("X62").Value = 39
$ws.Range("AB62").Value = 192
$ws.Range("X63").Value = 38
$ws.Range("AB63").Value = 184
$ws.Range("AB64").Value = 184
